# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" between "2021-Q3" and "总计" and fill
#    it with the fund-holding detail rows for the 2022-Q1 quarter (same
#    layout/styling as the "2021-Q3" sheet).
# 2) Insert a new summary row for "2022-Q1" at the top of the "总计" sheet's
#    data (pushing the existing 2021-Q3 / 2021-Q1 rows down by one row).

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing it to be stored as text
# (so numeric-looking strings like "3.07" or "001170" are not silently
# turned into real numbers), without leaving a lingering cell style behind.
function Set-TextValue($rng, $value) {
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------
# Step 1: create the new "2022-Q1" worksheet right after "2021-Q3"
# ---------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q3Sheet)
$newSheet.Name = "2022-Q1"

# Header row - same wording/format as sheet "2021-Q3", plus the two extra
# columns ("基金规模" and "仓位排名" already exist there, only the column
# order/content differs per the target data).
Set-TextValue $newSheet.Range("B1") "基金代码"
Set-TextValue $newSheet.Range("C1") "基金名称"
Set-TextValue $newSheet.Range("D1") "基金规模"
Set-TextValue $newSheet.Range("E1") "股票总仓位"
Set-TextValue $newSheet.Range("F1") "仓位占比"
Set-TextValue $newSheet.Range("G1") "持有市值(亿元)"
Set-TextValue $newSheet.Range("H1") "仓位排名"

# Copy the header/index-column formatting (bold, thin border, centered/top
# alignment) from the "2021-Q3" sheet so the new sheet matches its look.
$q3Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)

# Row 2
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "001170"
Set-TextValue $newSheet.Range("C2") "泰达宏利复兴伟业灵活配置混合"
Set-TextValue $newSheet.Range("D2") "3.07"
Set-TextValue $newSheet.Range("E2") "92.00"
Set-TextValue $newSheet.Range("F2") "4.04"
Set-TextValue $newSheet.Range("G2") "0.1240"
$newSheet.Range("H2").Value = 10

# Row 3
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "002289"
Set-TextValue $newSheet.Range("C3") "华商改革创新股票"
Set-TextValue $newSheet.Range("D3") "1.14"
Set-TextValue $newSheet.Range("E3") "90.69"
Set-TextValue $newSheet.Range("F3") "3.27"
Set-TextValue $newSheet.Range("G3") "0.0373"
$newSheet.Range("H3").Value = 7

# Row 4
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "001914"
Set-TextValue $newSheet.Range("C4") "中信建投聚利混合A"
Set-TextValue $newSheet.Range("D4") "0.13"
Set-TextValue $newSheet.Range("E4") "39.07"
Set-TextValue $newSheet.Range("F4") "2.65"
Set-TextValue $newSheet.Range("G4") "0.0034"
$newSheet.Range("H4").Value = 2

# Row 5
$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet.Range("B5") "000041"
Set-TextValue $newSheet.Range("C5") "华夏全球精选股票(QDII)"
Set-TextValue $newSheet.Range("D5") "0.02"
Set-TextValue $newSheet.Range("E5") "39.07"
Set-TextValue $newSheet.Range("F5") "2.65"
Set-TextValue $newSheet.Range("G5") "0.0005"
$newSheet.Range("H5").Value = 2

# Copy the index-column formatting (A column) for the data rows too.
$q3Sheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------
# Step 2: update the "总计" (summary) sheet with the new 2022-Q1 row
# ---------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Push the existing rows down by one: old row 3 (2021-Q1) -> row 4
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q1"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.01

# old row 2 (2021-Q3) -> row 3
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q3"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.11

# new row 2: 2022-Q1 summary
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.17

# A4 is a brand-new cell (the sheet previously only had rows 1-3), so copy
# the existing index-column style (from A2, which already carries it) onto
# the newly created A4 cell.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A4").PasteSpecial($xlPasteFormats)

# Restore the originally active sheet/selection
$wb.Worksheets.Item("2021-Q1").Activate()
